{"js": "// Change 1: \"A Q-Q Plot will help us know if the data are normally distributed.\"\n// becomes \"A histogram will help us know if the data are normally distributed.\"\n// (only in the item 10 paragraph that contains \"10 . A Q-Q Plot\")\nconst search1 = context.document.body.search(\"10 . A Q-Q Plot will help us know if the data are normally distributed.\", { matchCase: true });\nsearch1.load(\"text\");\nawait context.sync();\n\nif (search1.items.length > 0) {\n  search1.items[0].insertText(\n    \"10 . A histogram will help us know if the data are normally distributed.\",\n    \"Replace\"\n  );\n}\n\n// Change 2: expand/rewrite part \"b.\" explanation for non-certified/certified cars.\nconst search2 = context.document.body.search(\n  \"b. The sample size for non-certified cars is large, so we can assume that it is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution, but a Q-Q Plot will help us know more confidently.\",\n  { matchCase: true }\n);\nsearch2.load(\"text\");\nawait context.sync();\n\nif (search2.items.length > 0) {\n  search2.items[0].insertText(\n    \"b. The sample size for non-certified cars is large, so we can assume that the sampling distribution of the sample mean is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution for the mean, but a histogram of the data will help us know if the data is normally distributed. If it is, then the distribution of the sample mean will also be normally distributed.\",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: item 10 explanation - \"Q-Q Plot\" -> \"histogram\"\n$range1 = $d.Content\n$oldText1 = \"10 . A Q-Q Plot will help us know if the data are normally distributed.\"\n$newText1 = \"10 . A histogram will help us know if the data are normally distributed.\"\n$range1.Find.ClearFormatting()\n$range1.Find.Replacement.ClearFormatting()\n$range1.Find.Execute($oldText1, $false, $false, $false, $false, $false, $true, 1, $false, $newText1, 2)\n\n# Change 2: part b. explanation for non-certified/certified cars - rewritten/expanded\n$range2 = $d.Content\n$oldText2 = \"b. The sample size for non-certified cars is large, so we can assume that it is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution, but a Q-Q Plot will help us know more confidently.\"\n$newText2 = \"b. The sample size for non-certified cars is large, so we can assume that the sampling distribution of the sample mean is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution for the mean, but a histogram of the data will help us know if the data is normally distributed. If it is, then the distribution of the sample mean will also be normally distributed.\"\n$range2.Find.ClearFormatting()\n$range2.Find.Replacement.ClearFormatting()\n$range2.Find.Execute($oldText2, $false, $false, $false, $false, $false, $true, 1, $false, $newText2, 2)\n\n$d.Save()\n"}
